# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F values for rows 2-20
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 1112
$wsExpo.Range("F3").Value  = 412
$wsExpo.Range("F4").Value  = 248
$wsExpo.Range("F6").Value  = 12066
$wsExpo.Range("F7").Value  = 46
$wsExpo.Range("F8").Value  = 85
$wsExpo.Range("F9").Value  = 11819
$wsExpo.Range("F10").Value = 4756
$wsExpo.Range("F11").Value = 545
$wsExpo.Range("F12").Value = 73
$wsExpo.Range("F13").Value = 20
$wsExpo.Range("F14").Value = 420
$wsExpo.Range("F15").Value = 85
$wsExpo.Range("F19").Value = 56
$wsExpo.Range("F20").Value = 5217

# Sheet "全部类型" (All types) - column F values for rows 2-22
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value   = 1112
$wsAll.Range("F3").Value   = 412
$wsAll.Range("F4").Value   = 248
$wsAll.Range("F8").Value   = 12066
$wsAll.Range("F9").Value   = 46
$wsAll.Range("F10").Value  = 85
$wsAll.Range("F11").Value  = 11819
$wsAll.Range("F12").Value  = 4756
$wsAll.Range("F13").Value  = 545
$wsAll.Range("F14").Value  = 73
$wsAll.Range("F15").Value  = 20
$wsAll.Range("F16").Value  = 420
$wsAll.Range("F17").Value  = 85
$wsAll.Range("F21").Value  = 56
$wsAll.Range("F22").Value  = 5217

$wb.Save()
